$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 0.003173959397873089
$ws.Range("E2").Value = 0.01857558657652227
$ws.Range("G2").Value = 0.5
$ws.Range("H2").Value = 21658.18

# Row 3
$ws.Range("D3").Value = 0.001358097442220604
$ws.Range("E3").Value = 0.01592194719999851

# Row 4
$ws.Range("D4").Value = 0.001279591134608583
$ws.Range("E4").Value = 0.02234653965257309

# Row 5
$ws.Range("D5").Value = 0.0009345652867390518
$ws.Range("E5").Value = 0.02036751953979167

# Row 6
$ws.Range("D6").Value = 0.0021229050052257
$ws.Range("E6").Value = 0.02360997313468402

# Row 7
$ws.Range("D7").Value = 0.003779962272645701
$ws.Range("E7").Value = 0.0492683155047368
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0

# Row 8
$ws.Range("D8").Value = 0.0001458607358971155
$ws.Range("E8").Value = 0.01794364018502862

# Row 9
$ws.Range("D9").Value = 0.004854294480505244
$ws.Range("E9").Value = 0.0423008279237936
$ws.Range("G9").Value = 0.06
$ws.Range("H9").Value = 2598.98

# Row 10
$ws.Range("D10").Value = 0.009271276528460179
$ws.Range("E10").Value = 0.06758358167121899
$ws.Range("G10").Value = 0.13
$ws.Range("H10").Value = 5631.13

# Row 11
$ws.Range("D11").Value = 0.004078503750362078
$ws.Range("E11").Value = 0.02793348604957986
$ws.Range("G11").Value = 0.31
$ws.Range("H11").Value = 13428.07
